$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + report week dates) ---
# "Volume 30   Number  48" -> "...49"
$ws.Range("A8").Characters(21, 2).Text = "49"
# "Report Covering the Week  11/27/2023  Through  12/3/2023"
# -> "...12/4/2023  Through  12/10/2023"
# (replace the right-hand date first so the left-hand offset stays valid)
$ws.Range("C9").Characters(48, 9).Text = "12/10/2023"
$ws.Range("C9").Characters(27, 10).Text = "12/4/2023"

# --- Cells reverting to a "blank" (0) / "undefined" (***.*) text placeholder, style 14 ---
# Copy(destination) clones value+format together from an existing placeholder cell,
# so the literal text is preserved (a plain .Value="0" would be re-coerced to a number).
$ws.Range("C14").Copy($ws.Range("C19"))
$ws.Range("C14").Copy($ws.Range("G19"))
$ws.Range("E14").Copy($ws.Range("H19"))
$ws.Range("C14").Copy($ws.Range("D24"))
$ws.Range("E14").Copy($ws.Range("E24"))

# --- Cells changing number format (same font group, different numFmt) ---
$ws.Range("C17").NumberFormat = "#,##0"
$ws.Range("C17").Value = 1
$ws.Range("D17").NumberFormat = "#,##0"
$ws.Range("D17").Value = 1
$ws.Range("E17").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E17").Value = 0
$ws.Range("F17").NumberFormat = "#,##0"
$ws.Range("F17").Value = 1
$ws.Range("G17").NumberFormat = "#,##0"
$ws.Range("G17").Value = 1
$ws.Range("H17").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("H17").Value = 0
$ws.Range("D18").NumberFormat = "#,##0"
$ws.Range("D18").Value = 2
$ws.Range("E18").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E18").Value = -100
$ws.Range("G18").NumberFormat = "#,##0"
$ws.Range("G18").Value = 2
$ws.Range("H18").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("H18").Value = -100
$ws.Range("D21").NumberFormat = "#,##0"
$ws.Range("D21").Value = 3
$ws.Range("E21").NumberFormat = "#,##0.00;`"-`"#,##0.00"
$ws.Range("E21").Value = -66.666666666666
$ws.Range("C24").NumberFormat = "#,##0"
$ws.Range("C24").Value = 1
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("D30").Value = 1
$ws.Range("E30").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E30").Value = -100
$ws.Range("G30").NumberFormat = "#,##0"
$ws.Range("G30").Value = 1
$ws.Range("H30").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("H30").Value = 0
$ws.Range("J30").NumberFormat = "#,##0"
$ws.Range("J30").Value = 1
$ws.Range("K30").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("K30").Value = 300

# --- Cells with unchanged style, value update only ---
$ws.Range("I17").Value = 9
$ws.Range("J17").Value = 14
$ws.Range("K17").Value = -35.714285714285
$ws.Range("L17").Value = -10
$ws.Range("M17").Value = 80
$ws.Range("N17").Value = -75.675675675675
$ws.Range("J18").Value = 4
$ws.Range("K18").Value = 25
$ws.Range("L18").Value = 150
$ws.Range("L19").Value = 84.615384615384
$ws.Range("N19").Value = -72.413793103448
$ws.Range("C21").Value = 1
$ws.Range("F21").Value = 6
$ws.Range("G21").Value = 3
$ws.Range("H21").Value = 100
$ws.Range("I21").Value = 84
$ws.Range("J21").Value = 75
$ws.Range("K21").Value = 12
$ws.Range("L21").Value = 44.827586206896
$ws.Range("M21").Value = -14.285714285714
$ws.Range("N21").Value = -81.778741865509
$ws.Range("F24").Value = 2
$ws.Range("H24").Value = 100
$ws.Range("I24").Value = 40
$ws.Range("K24").Value = 33.333333333333
$ws.Range("L24").Value = 21.212121212121
$ws.Range("M24").Value = -60.39603960396
$ws.Range("G25").Value = 3
$ws.Range("H25").Value = -33.333333333333
$ws.Range("J25").Value = 31
$ws.Range("K25").Value = 51.612903225806
$ws.Range("L27").Value = 200
$ws.Range("L30").Value = 33.333333333333
